$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "62.542.76"
$r.Style = "Normal"
$ws.Range("E2").Value = "  +0.55%  "

# Row 3
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "2.441.17"
$r.Style = "Normal"
$ws.Range("E3").Value = "  +1.03%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "566.72"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "

# Row 6
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "145.54"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +1.85%  "

# Row 7
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "

# Row 8
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.534"
$r.Style = "Normal"
$ws.Range("E8").Value = "  +0.66%  "

# Row 9
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.112"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +2.16%  "

# Row 10
$ws.Range("E10").Value = "  +0.29%  "

# Row 11
$ws.Range("E11").Value = "  -1.23%  "

# Row 12
$ws.Range("E12").Value = "  -0.09%  "

# Row 13
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "26.85"
$r.Style = "Normal"
$ws.Range("E13").Value = "  +4.65%  "

# Row 14
$ws.Range("E14").Value = "  +4.67%  "

# Row 15
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "2.833.81"
$r.Style = "Normal"
$ws.Range("E15").Value = "  -0.65%  "

# Row 16
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "62.513.06"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "

# Row 17
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "2.445.11"
$r.Style = "Normal"
$ws.Range("E17").Value = "  +1.36%  "

# Row 18
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "11.27"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -0.42%  "

# Row 19
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "6.96"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +1.64%  "

# Row 20
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "324.31"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "

# Row 21
$ws.Range("E21").Value = "  -0.07%  "

# Row 22
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "67.43"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +2.52%  "

# Row 24
$ws.Range("E24").Value = "  +2.26%  "

# Row 25
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "8.76"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -2.84%  "

# Row 26
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "0.0₃0987"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +4.20%  "

# Row 27
$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "561.15"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -2.87%  "

# Row 28
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "2.555.21"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +0.95%  "

# Row 29
$ws.Range("E29").Value = "  -0.18%  "

# Row 30
$ws.Range("E30").Value = "  +1.50%  "

# Row 31
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "1.46"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +1.87%  "

# Row 32
$ws.Range("E32").Value = "  -0.81%  "

# Row 33
$ws.Range("E33").Value = "  +0.59%  "

# Row 34
$ws.Range("E34").Value = "  +0.89%  "

# Row 35
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "4.89"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +3.55%  "

# Row 36
$ws.Range("E36").Value = "  -0.15%  "

# Row 37
$ws.Range("E37").Value = "  +0.37%  "

# Row 38
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "5.48"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "

# Row 39
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "18.80"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "

# Row 40
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "150.35"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -1.11%  "

# Row 41
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "1.82"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +0.50%  "

# Row 42
$ws.Range("E42").Value = "  +0.05%  "

# Row 43
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "2.41"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +5.13%  "

# Row 44
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "148.94"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +0.29%  "

# Row 45
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "3.69"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +1.02%  "

# Row 46
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.0537"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "

# Row 47
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "20.47"
$r.Style = "Normal"
$ws.Range("E47").Value = "  +2.11%  "

# Row 48
$ws.Range("E48").Value = "  +1.36%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.0232"
$r.Style = "Normal"
$ws.Range("E49").Value = "  +1.90%  "

# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.0928"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "

# Row 51
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "11.59"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +0.41%  "

Write-Output "Applied all changes"